# LBS-stories-LP1.xlsx : "get carte + paiement carte"
#
# Adds the new API routes for the "carte de fidelite" (loyalty card)
# endpoints (GET auth, GET carte, POST paiement) to the user-story
# table, marks rows 25-27 as "not done" (Realise = 0), extends the
# points-total formula in E44 to cover the newly completed rows, and
# adds a new Print_Area_0_0_0_0_0 defined name (mirroring the existing
# Print_Area_* duplicates already present in the workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 25-27 : mark "Realise" (H) explicitly as 0 (not yet done) ---
$ws.Range("H25").Value = 0
$ws.Range("H26").Value = 0
$ws.Range("H27").Value = 0

# --- Row 30 : "s'authentifier pour utiliser sa carte" -> GET /cartes/{id}/auth ---
$ws.Range("F30").Value = "GET"
$ws.Range("G30").Value = "/cartes/{id}/auth"

# --- Row 31 : "acceder a sa carte de fidelite" -> GET /cartes/{id} ---
$ws.Range("F31").Value = "GET"
$ws.Range("G31").Value = "/cartes/{id}"
$ws.Range("I31").Value = "Baptiste/Islam/Mohammed"

# --- Row 32 : "paiement fidelise" -> POST /cartes/{id}/paiement, done, authors ---
$ws.Range("F32").Value = "POST"
$ws.Range("G32").Value = "/cartes/{id}/paiement"
$ws.Range("H32").Value = 1
$ws.Range("I32").Value = "Baptiste/Islam"

# --- Extend the points-total formula to include the newly scored rows ---
$ws.Range("E44").Formula = "=E7*H7+E8*H8+E9*H9+E10*H10+E11*H11+E12*H12+E13*H13+E14*H14+E15*H15+E16*H16+E17*H17+E18*H18+E19*H19+E20*H20+E21*H21+E22*H22+E25*H25+E26*H26+E27*H27+E30*H30+E31*H31+E32*H32+E33*H33+E36*H36+E37*H37+E38*H38+E39*H39+E40*H40+E41*H41"

# --- Add the new duplicated print-area defined name (sheet scoped) ---
$ws.Names.Add("_xlnm.Print_Area_0_0_0_0_0", "=Feuille1!`$A`$1:`$G`$47")
